# Insert a new worksheet "Prefabs View" between "Animals" and "DayNightCycle",
# matching the commit "feat: Entity view and view classes".

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Prefabs View" sheet, positioned right before DayNightCycle ---
$dayNight = $wb.Worksheets.Item("DayNightCycle")
$prefabsView = $wb.Worksheets.Add($dayNight)
$prefabsView.Name = "Prefabs View"

# Header row
$prefabsView.Cells.Item(1, 1).Value = "IDS"
$prefabsView.Cells.Item(1, 2).Value = "Architecture ID"
$prefabsView.Cells.Item(1, 3).Value = "Prefab resource path"

# Data row: the Monkey entity view
$prefabsView.Cells.Item(2, 1).Value = "Monkey view"
$prefabsView.Cells.Item(2, 2).Value = "Monkey"
$prefabsView.Cells.Item(2, 3).Value = "Entities/LivingEntities/Animals/Monkey.prefab"

# --- 2. Tile Types sheet picked up a custom width on column D ---
$tileTypes = $wb.Worksheets.Item("Tile Types")
$tileTypes.Columns.Item(4).ColumnWidth = 14.25

# --- 3. Restore the original active sheet ---
$wb.Worksheets.Item("Animals").Activate()
